$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.992.39'
$ws.Range('E2').Value = '  -3.75%  '
$ws.Range('D3').Value = '1.871.17'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4358'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3757'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07486'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9355'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.80%  '
$ws.Range('E11').Value = '  -5.29%  '
$ws.Range('D12').Value = '1.916.12'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('E13').Value = '  -3.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.442'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06865'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '81.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009008'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.88'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.36%  '
$ws.Range('D21').Value = '27.987.93'
$ws.Range('E21').Value = '  -3.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').Value = '2.130.85'
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.036'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.74'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.613'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.708'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09024'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8139'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.818'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.184'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.963'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.13%  '
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  -2.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05515'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01979'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.958'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5274'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1701'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.009'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.776'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06756'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4882'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.59'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.37%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.914'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -14.43%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.677'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.98%  '
$ws.Range('E51').Value = '  -0.11%  '
